# Macroferia Regional de Talca - Acelga: insert one new weekly price record.
#
# The sheet is a flat daily/weekly price log (rows 2..275, one row per
# observation). This edit inserts a brand-new observation at row 200,
# pushing every existing row from 200..275 down by one (to 201..276), and
# populates the new row 200 with its own date/volume/price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 200:275 down to 201:276, opening up a blank row 200.
$ws.Rows("200:200").Insert()

# The new row 200 shares every descriptive column with the rest of this
# market/category block; only the date, volume and price columns differ.
$ws.Cells.Item(200, 1).Value2  = $ws.Cells.Item(201, 1).Value2    # Mercado ID
$ws.Cells.Item(200, 2).Value2  = $ws.Cells.Item(201, 2).Value2    # Mercado
$ws.Cells.Item(200, 3).Value2  = $ws.Cells.Item(201, 3).Value2    # Región
$ws.Cells.Item(200, 4).Value2  = 44755                             # Fecha
$ws.Cells.Item(200, 5).Value2  = $ws.Cells.Item(201, 5).Value2    # Codreg
$ws.Cells.Item(200, 6).Value2  = $ws.Cells.Item(201, 6).Value2    # Categoría ID
$ws.Cells.Item(200, 7).Value2  = $ws.Cells.Item(201, 7).Value2    # Categoría
$ws.Cells.Item(200, 8).Value2  = $ws.Cells.Item(201, 8).Value2    # Variedad
$ws.Cells.Item(200, 9).Value2  = $ws.Cells.Item(201, 9).Value2    # Calidad
$ws.Cells.Item(200, 10).Value2 = 300                                # Volumen
$ws.Cells.Item(200, 11).Value2 = 3000                                # Precio mínimo
$ws.Cells.Item(200, 12).Value2 = 3000                                # Precio máximo
$ws.Cells.Item(200, 13).Value2 = 3000                                # Precio promedio ponderado
$ws.Cells.Item(200, 14).Value2 = $ws.Cells.Item(201, 14).Value2   # Unidad de comercialización
$ws.Cells.Item(200, 15).Value2 = $ws.Cells.Item(201, 15).Value2   # Origen
$ws.Cells.Item(200, 16).Value2 = 750                                 # Precio $/Kg
$ws.Cells.Item(200, 17).Value2 = $ws.Cells.Item(201, 17).Value2   # Kg o Unidades
$ws.Cells.Item(200, 18).Value2 = $ws.Cells.Item(201, 18).Value2   # Clasificación
